# Update the "想去人数" (want-to-go count) figures for three camp/expo
# events that changed between the previous and current data pull.
#
# Each event appears on two sheets: "展览" (the exhibitions-only view) and
# "全部类型" (the all-types view), so both copies need updating.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 7693   # was 7692
    $ws.Range("F5").Value = 473    # was 472
    $ws.Range("F6").Value = 4411   # was 4400
}
